# Update the division problems in the practice table.
# The table has 5 "data" rows (1, 5, 9, 13, 17), each with 5 columns,
# interleaved with empty spacer rows. We update each cell's text in
# place so that formatting (rFonts/sz) carried by the existing run is
# preserved.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "73÷4="; New = "50÷4=" },
    @{ Row = 1;  Col = 2; Old = "62÷7="; New = "13÷6=" },
    @{ Row = 1;  Col = 3; Old = "33÷3="; New = "38÷9=" },
    @{ Row = 1;  Col = 4; Old = "33÷3="; New = "93÷5=" },
    @{ Row = 1;  Col = 5; Old = "30÷2="; New = "97÷7=" },

    @{ Row = 5;  Col = 1; Old = "13÷8="; New = "91÷5=" },
    @{ Row = 5;  Col = 2; Old = "29÷3="; New = "40÷8=" },
    @{ Row = 5;  Col = 3; Old = "37÷6="; New = "55÷7=" },
    @{ Row = 5;  Col = 4; Old = "60÷5="; New = "27÷2=" },
    @{ Row = 5;  Col = 5; Old = "72÷3="; New = "16÷4=" },

    @{ Row = 9;  Col = 1; Old = "27÷4="; New = "92÷6=" },
    @{ Row = 9;  Col = 2; Old = "33÷4="; New = "24÷9=" },
    @{ Row = 9;  Col = 3; Old = "62÷8="; New = "61÷2=" },
    @{ Row = 9;  Col = 4; Old = "14÷9="; New = "47÷2=" },
    @{ Row = 9;  Col = 5; Old = "85÷4="; New = "69÷5=" },

    @{ Row = 13; Col = 1; Old = "56÷6="; New = "81÷2=" },
    @{ Row = 13; Col = 2; Old = "66÷8="; New = "42÷6=" },
    @{ Row = 13; Col = 3; Old = "79÷7="; New = "23÷8=" },
    @{ Row = 13; Col = 4; Old = "55÷3="; New = "72÷6=" },
    @{ Row = 13; Col = 5; Old = "86÷4="; New = "55÷8=" },

    @{ Row = 17; Col = 1; Old = "65÷4="; New = "31÷3=" },
    @{ Row = 17; Col = 2; Old = "81÷7="; New = "13÷7=" },
    @{ Row = 17; Col = 3; Old = "63÷2="; New = "35÷5=" },
    @{ Row = 17; Col = 4; Old = "78÷9="; New = "76÷9=" },
    @{ Row = 17; Col = 5; Old = "55÷7="; New = "52÷2=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $rng = $cell.Range
    # Cell.Range.Text includes the trailing cell-mark (CR + BEL); trim it
    # off so we only compare/replace the visible problem text.
    $current = $rng.Text
    $current = $current.TrimEnd([char]13, [char]7)
    if ($current -ne $u.Old) {
        Write-Host "WARNING: Row $($u.Row) Col $($u.Col) expected '$($u.Old)' but found '$current'"
    }
    $rng.Text = $u.New
}
